$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44260.50046984621
$ws.Range("D16:D29").Value = 44260.47913657407
$ws.Range("D30:D43").Value = 44260.45781210648
